$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Populate new columns F:I with RDS-file diet/virus counts ---

# Block 1 (row1 "TOTH" header, rows 3-5 DESeq2/EdgeR/Limma) - RDS "diet" files
$ws.Range("F3").Value = "RDC"
$ws.Range("G3").Value = "RDR"
$ws.Range("H3").Value = "RD_DIET"
$ws.Range("I3").Value = "Done"

$ws.Range("F4").Value = "REC"
$ws.Range("G4").Value = "RER"
$ws.Range("H4").Value = "RE_DIET"
$ws.Range("I4").Value = "Done"

$ws.Range("F5").Value = "RLC"
$ws.Range("G5").Value = "RLR"
$ws.Range("H5").Value = "RL_DIET"
$ws.Range("I5").Value = "Done"

# Block 2 (row7 "GALBRAITH" header, rows 8-10 DESeq2/EdgeR/Limma) - RDS "virus" files
$ws.Range("F8").Value = "RDV"
$ws.Range("G8").Value = "RDN"
$ws.Range("H8").Value = "RD_VIRUS"
$ws.Range("I8").Value = "Done"

$ws.Range("F9").Value = "REV"
$ws.Range("G9").Value = "REN"
$ws.Range("H9").Value = "RE_VIRUS"
$ws.Range("I9").Value = "Done"

$ws.Range("F10").Value = "RLV"
$ws.Range("G10").Value = "RLN"
$ws.Range("H10").Value = "RL_VIRUS"
$ws.Range("I10").Value = "Done"

# Block 3 (row12 "TOTH PAIRS" header, rows 13-15 DESeq2/EdgeR/Limma) - Galbraith RDS files
$ws.Range("F13").Value = "GDV"
$ws.Range("G13").Value = "GDC"
$ws.Range("H13").Value = "GD_VIRUS"

$ws.Range("F14").Value = "GEV"
$ws.Range("G14").Value = "GEC"
$ws.Range("H14").Value = "GE_VIRUS"

$ws.Range("F15").Value = "GLV"
$ws.Range("G15").Value = "GLC"
$ws.Range("H15").Value = "GL_VIRUS"

# --- View/selection bookkeeping ---
$ws.Range("I8").Select()
